$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Stock" row's Caso 1 value from "xsd:int" to "xsd:integer"
$ws.Range("B7").Value = "xsd:integer"

# Update the active cell selection from C15 to C14
$ws.Range("C14").Select()
